$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 and 4 (keep only header + 1 data row)
$ws.Rows("3:4").Delete()

# Update row 2 values
$ws.Range("A2").Value = "tk2504200@gmail.com"
$ws.Range("B2").Value = "Khanh2504"
$ws.Range("C2").Value = "35"
$ws.Range("D2").Value = "73"
$ws.Range("E2").Value = "222"
$ws.Range("F2").Value = "1285"
$ws.Range("G2").Value = "1"

$ws.Range("F2").Select()

$ws.Columns("A").ColumnWidth = 24.21875
